$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Match formatting used by the rest of the table before assigning values, so
# no ad-hoc number-format style gets minted for the new date cell: date cell
# + day marker use style 4 (custom date number format), notes column uses
# style 5 (wrap text).
$ws.Range("A42").NumberFormat = $ws.Range("A41").NumberFormat
$ws.Range("B42").NumberFormat = $ws.Range("B41").NumberFormat
$ws.Range("E42").WrapText = $true

# New entry: Tue July 1 2024 ("M" day-of-week marker), 6 hours, note about
# cleaning the professors table and merging it into the main students dataframe.
$ws.Range("A42").Value = Get-Date -Year 2024 -Month 7 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("B42").Value = "M"
$ws.Range("C42").Value = 6
$ws.Range("E42").Value = "cleaning professors table, merging professors into main students dataframe"

$ws.Rows.Item(42).RowHeight = 28.5

$ws.Range("E42").Select()
